$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 46.62552400000001
$ws.Range("H2").Value = 139.876572
$ws.Range("I2").Value = 0.1981356926336952
$ws.Range("J2").Value = 0.1981356926336952
$ws.Range("M2").Value = 2.326231
$ws.Range("N2").Value = 6.978693
$ws.Range("O2").Value = 0.0902108727640205
$ws.Range("P2").Value = 0.09021087276402051
$ws.Range("Q2").Value = 108.461739320044
$ws.Range("R2").Value = 976.155653880396
$ws.Range("S2").Value = 0.01787399375818935
$ws.Range("T2").Value = 0.01787399375818935

# Row 3
$ws.Range("G3").Value = 46.62552400000001
$ws.Range("H3").Value = 139.876572
$ws.Range("I3").Value = 0.1981356926336952
$ws.Range("J3").Value = 0.1981356926336952
$ws.Range("O3").Value = 0.4488373460538801
$ws.Range("P3").Value = 0.4488373460538801
$ws.Range("Q3").Value = 539.6431464768227
$ws.Range("R3").Value = 4856.788318291404
$ws.Range("S3").Value = 0.08893069844025508
$ws.Range("T3").Value = 0.08893069844025507

# Row 4
$ws.Range("G4").Value = 46.62552400000001
$ws.Range("H4").Value = 139.876572
$ws.Range("I4").Value = 0.1981356926336952
$ws.Range("J4").Value = 0.1981356926336952
$ws.Range("M4").Value = 8.667390666666668
$ws.Range("N4").Value = 26.002172
$ws.Range("O4").Value = 0.3361200485363344
$ws.Range("P4").Value = 0.3361200485363344
$ws.Range("Q4").Value = 404.1216315460427
$ws.Range("R4").Value = 3637.094683914384
$ws.Range("S4").Value = 0.06659737862481788
$ws.Range("T4").Value = 0.06659737862481785

# Row 5
$ws.Range("G5").Value = 46.62552400000001
$ws.Range("H5").Value = 139.876572
$ws.Range("I5").Value = 0.1981356926336952
$ws.Range("J5").Value = 0.1981356926336952
$ws.Range("M5").Value = 3.218985
$ws.Range("N5").Value = 9.656955
$ws.Range("O5").Value = 0.124831732645765
$ws.Range("P5").Value = 0.124831732645765
$ws.Range("Q5").Value = 150.08686237314
$ws.Range("R5").Value = 1350.78176135826
$ws.Range("S5").Value = 0.02473362181043291
$ws.Range("T5").Value = 0.02473362181043291

# Row 6
$ws.Range("I6").Value = 0.441606206212991
$ws.Range("J6").Value = 0.4416062062129909
$ws.Range("M6").Value = 2.326231
$ws.Range("N6").Value = 6.978693
$ws.Range("O6").Value = 0.0902108727640205
$ws.Range("P6").Value = 0.09021087276402051
$ws.Range("Q6").Value = 241.7402769976313
$ws.Range("R6").Value = 2175.662492978682
$ws.Range("S6").Value = 0.03983768128048193
$ws.Range("T6").Value = 0.03983768128048193

# Row 7
$ws.Range("I7").Value = 0.441606206212991
$ws.Range("J7").Value = 0.4416062062129909
$ws.Range("O7").Value = 0.4488373460538801
$ws.Range("P7").Value = 0.4488373460538801
$ws.Range("Q7").Value = 1202.760388382157
$ws.Range("S7").Value = 0.1982093575975614
$ws.Range("T7").Value = 0.1982093575975614

# Row 8
$ws.Range("I8").Value = 0.441606206212991
$ws.Range("J8").Value = 0.4416062062129909
$ws.Range("M8").Value = 8.667390666666668
$ws.Range("N8").Value = 26.002172
$ws.Range("O8").Value = 0.3361200485363344
$ws.Range("P8").Value = 0.3361200485363344
$ws.Range("Q8").Value = 900.7090957891477
$ws.Range("R8").Value = 8106.381862102329
$ws.Range("S8").Value = 0.148432699466257
$ws.Range("T8").Value = 0.148432699466257

# Row 9
$ws.Range("I9").Value = 0.441606206212991
$ws.Range("J9").Value = 0.4416062062129909
$ws.Range("M9").Value = 3.218985
$ws.Range("N9").Value = 9.656955
$ws.Range("O9").Value = 0.124831732645765
$ws.Range("P9").Value = 0.124831732645765
$ws.Range("Q9").Value = 334.51464001263
$ws.Range("R9").Value = 3010.63176011367
$ws.Range("S9").Value = 0.05512646786869065
$ws.Range("T9").Value = 0.05512646786869065

# Row 10
$ws.Range("G10").Value = 49.27528633333333
$ws.Range("H10").Value = 147.825859
$ws.Range("I10").Value = 0.2093958876983056
$ws.Range("J10").Value = 0.2093958876983056
$ws.Range("M10").Value = 2.326231
$ws.Range("N10").Value = 6.978693
$ws.Range("O10").Value = 0.0902108727640205
$ws.Range("P10").Value = 0.09021087276402051
$ws.Range("Q10").Value = 114.6256986024763
$ws.Range("R10").Value = 1031.631287422287
$ws.Range("S10").Value = 0.01888978578246097
$ws.Range("T10").Value = 0.01888978578246097

# Row 11
$ws.Range("G11").Value = 49.27528633333333
$ws.Range("H11").Value = 147.825859
$ws.Range("I11").Value = 0.2093958876983056
$ws.Range("J11").Value = 0.2093958876983056
$ws.Range("O11").Value = 0.4488373460538801
$ws.Range("P11").Value = 0.4488373460538801
$ws.Range("Q11").Value = 570.3114577428958
$ws.Range("R11").Value = 5132.803119686062
$ws.Range("S11").Value = 0.09398469450910382
$ws.Range("T11").Value = 0.09398469450910381

# Row 12
$ws.Range("G12").Value = 49.27528633333333
$ws.Range("H12").Value = 147.825859
$ws.Range("I12").Value = 0.2093958876983056
$ws.Range("J12").Value = 0.2093958876983056
$ws.Range("M12").Value = 8.667390666666668
$ws.Range("N12").Value = 26.002172
$ws.Range("O12").Value = 0.3361200485363344
$ws.Range("P12").Value = 0.3361200485363344
$ws.Range("Q12").Value = 427.0881568628609
$ws.Range("R12").Value = 3843.793411765748
$ws.Range("S12").Value = 0.07038215593646331
$ws.Range("T12").Value = 0.0703821559364633

# Row 13
$ws.Range("G13").Value = 49.27528633333333
$ws.Range("H13").Value = 147.825859
$ws.Range("I13").Value = 0.2093958876983056
$ws.Range("J13").Value = 0.2093958876983056
$ws.Range("M13").Value = 3.218985
$ws.Range("N13").Value = 9.656955
$ws.Range("O13").Value = 0.124831732645765
$ws.Range("P13").Value = 0.124831732645765
$ws.Range("Q13").Value = 158.616407577705
$ws.Range("R13").Value = 1427.547668199345
$ws.Range("S13").Value = 0.02613925147027752
$ws.Range("T13").Value = 0.02613925147027751

# Row 14
$ws.Range("G14").Value = 35.50107333333333
$ws.Range("H14").Value = 106.50322
$ws.Range("I14").Value = 0.1508622134550082
$ws.Range("J14").Value = 0.1508622134550081
$ws.Range("M14").Value = 2.326231
$ws.Range("N14").Value = 6.978693
$ws.Range("O14").Value = 0.0902108727640205
$ws.Range("P14").Value = 0.09021087276402051
$ws.Range("Q14").Value = 82.58369732127332
$ws.Range("R14").Value = 743.25327589146
$ws.Range("S14").Value = 0.01360941194288824
$ws.Range("T14").Value = 0.01360941194288824

# Row 15
$ws.Range("G15").Value = 35.50107333333333
$ws.Range("H15").Value = 106.50322
$ws.Range("I15").Value = 0.1508622134550082
$ws.Range("J15").Value = 0.1508622134550081
$ws.Range("O15").Value = 0.4488373460538801
$ws.Range("P15").Value = 0.4488373460538801
$ws.Range("Q15").Value = 410.8889139112821
$ws.Range("R15").Value = 3698.00022520154
$ws.Range("S15").Value = 0.06771259550695982
$ws.Range("T15").Value = 0.06771259550695982

# Row 16
$ws.Range("G16").Value = 35.50107333333333
$ws.Range("H16").Value = 106.50322
$ws.Range("I16").Value = 0.1508622134550082
$ws.Range("J16").Value = 0.1508622134550081
$ws.Range("M16").Value = 8.667390666666668
$ws.Range("N16").Value = 26.002172
$ws.Range("O16").Value = 0.3361200485363344
$ws.Range("P16").Value = 0.3361200485363344
$ws.Range("Q16").Value = 307.7016716659822
$ws.Range("R16").Value = 2769.31504499384
$ws.Range("S16").Value = 0.05070781450879618
$ws.Range("T16").Value = 0.05070781450879617

# Row 17
$ws.Range("G17").Value = 35.50107333333333
$ws.Range("H17").Value = 106.50322
$ws.Range("I17").Value = 0.1508622134550082
$ws.Range("J17").Value = 0.1508622134550081
$ws.Range("M17").Value = 3.218985
$ws.Range("N17").Value = 9.656955
$ws.Range("O17").Value = 0.124831732645765
$ws.Range("P17").Value = 0.124831732645765
$ws.Range("Q17").Value = 114.2774225439
$ws.Range("R17").Value = 1028.4968028951
$ws.Range("S17").Value = 0.01883239149636391
$ws.Range("T17").Value = 0.01883239149636391
